$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "A" column (column B) entirely - it is being dropped from the table
$ws.Range("B1:B4").EntireColumn.Delete() | Out-Null

# Remove the "A Lag" row (now row 2 after the column delete) entirely
$ws.Range("A2:C2").EntireRow.Delete() | Out-Null

# Update the remaining data values (FFR Lag / LF Lag rows) with the new figures
$ws.Range("B2").Value = "1.88***"
$ws.Range("B3").Value = "3.77*"
$ws.Range("C2").Value = "0.47***"
$ws.Range("C3").Value = "0.76*"
